# Refresh the "Updated symbol list" crypto snapshot (cryptos.xlsx / Sheet1)
# with the GitHub Actions run from Fri Feb 17 14:40:00 UTC 2023:
#   - Price (col D) and Volume(1h) (col E) refreshed for most rows.
#   - Rows 11/12 swap which coin (BitrueCoin / MandalaExchangeToken) is
#     listed first, with fresh Coin/Link/Price/Volume data for each.
#
# The source cells are plain text (<is><t>...</t></is>, no numFmt) even
# though most of them look numeric ("310.71", "-1.88%", "0.2020", ...).
# Typing a numeric-looking string into a General-formatted cell makes
# Excel auto-convert it to a real number, which would both round the
# value (binary float error) and silently drop meaningful trailing
# zeros (e.g. "0.2020" -> 0.202). To keep every digit exactly as in the
# source data, each such cell is switched to Text format ("@") right
# before its value is written so Excel stores/echoes it verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.88%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '54.06'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '12.30%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.098'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.41%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07861'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.69%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.507'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.88%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.363'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.95%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.588'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.11%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1237'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.09%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2020'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.72%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09514'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.10%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04714'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.90%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1044'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.25%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001268'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-4.20%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005760'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.35%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,008.74%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.338'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.33%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.415'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.55%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3448'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.04%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.964'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.67%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1367'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.03%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3074'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.97%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04161'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.14%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001255'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-5.05%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003924'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.28%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001343'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.74%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02616'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-1.49%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05965'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.23%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01048'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.75%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007924'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.22%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1425'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.54%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008194'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.38%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008425'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.71%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3126'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.99%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007275'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.84%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000746'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.89%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05642'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '2.68%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002606'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-35.06%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002089'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.89%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001989'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.89%'

Write-Output "Updated symbol list: 81 cells refreshed on Sheet1."
